$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "43.819.99"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +4.71%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.261.09"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +1.46%  "
$ws.Range("E4").Value = "  -0.01%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "230.47"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.59%  "
$ws.Range("E6").Value = "  -0.10%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "61.20"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -0.68%  "
$ws.Range("E8").Value = "  +0.03%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.419"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +4.14%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "57.95"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -2.19%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0931"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +3.42%  "
$ws.Range("E12").Value = "  +0.23%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "2.599.43"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +1.43%  "
$ws.Range("E14").Value = "  -0.65%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "23.48"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +6.67%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "5.79"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +3.63%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.806"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +0.43%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "2.263.88"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +0.83%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "43.628.32"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +4.57%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.0₃0931"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +3.15%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "72.85"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +1.04%  "
$ws.Range("E22").Value = "  +2.82%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "254.48"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +1.75%  "
$ws.Range("E24").Value = "  -0.08%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.54"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +6.07%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.37"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +2.10%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "9.83"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +1.56%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "170.40"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +2.01%  "
$ws.Range("E29").Value = "  -1.04%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "20.45"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +2.30%  "
$ws.Range("E31").Value = "  +1.18%  "
$ws.Range("E32").Value = "  +0.40%  "
$ws.Range("E33").Value = "  +0.07%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "5.09"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +1.79%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "4.75"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +1.31%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.0657"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +3.15%  "
$ws.Range("E37").Value = "  -2.83%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.38"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +0.51%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "3.60"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -1.57%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.0249"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +4.04%  "
$ws.Range("E41").Value = "  -0.02%  "
$ws.Range("E42").Value = "  +1.39%  "
$ws.Range("E43").Value = "  -10.32%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.0981"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +0.60%  "
$ws.Range("E45").Value = "  -6.90%  "
$ws.Range("E46").Value = "  -1.11%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "97.94"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -1.03%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.472.11"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -0.66%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "16.61"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +0.65%  "
$ws.Range("E50").Value = "  +0.33%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "2.24"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +6.59%  "
